# Update cryptos list with latest price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.035.12'
$ws.Range("E2").Value = '  -0.27%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.550.36'
$ws.Range("E3").Value = '  -0.16%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.35'
$ws.Range("E5").Value = '  +1.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.64'
$ws.Range("E6").Value = '  +4.46%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  -0.76%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.73'
$ws.Range("E10").Value = '  +2.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0824'
$ws.Range("E11").Value = '  +1.53%  '

# Row 12
$ws.Range("E12").Value = '  +5.49%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.61'
$ws.Range("E13").Value = '  -2.01%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.945.29'
$ws.Range("E14").Value = '  +0.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.575.45'
$ws.Range("E15").Value = '  +1.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.86'
$ws.Range("E16").Value = '  +4.80%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.878'
$ws.Range("E17").Value = '  +0.43%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.119.93'
$ws.Range("E18").Value = '  -0.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("E19").Value = '  +4.18%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0990'
$ws.Range("E20").Value = '  +0.81%  '

# Row 21
$ws.Range("E21").Value = '  -0.45%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.01'
$ws.Range("E22").Value = '  -0.53%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.31'
$ws.Range("E23").Value = '  -2.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").Value = '  +1.39%  '

# Row 25
$ws.Range("E25").Value = '  -2.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.06'
$ws.Range("E26").Value = '  -6.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.15'
$ws.Range("E28").Value = '  +1.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.95'
$ws.Range("E29").Value = '  +1.43%  '

# Row 30
$ws.Range("E30").Value = '  -0.70%  '

# Row 31
$ws.Range("E31").Value = '  +0.30%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.69'
$ws.Range("E32").Value = '  +2.91%  '

# Row 33
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  -0.43%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.16'
$ws.Range("E34").Value = '  -0.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0806'
$ws.Range("E35").Value = '  +0.80%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.30'
$ws.Range("E36").Value = '  -2.82%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.84'
$ws.Range("E37").Value = '  +12.19%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.18'
$ws.Range("E38").Value = '  +11.62%  '

# Row 39
$ws.Range("E39").Value = '  -1.39%  '

# Row 40
$ws.Range("E40").Value = '  -0.33%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.13'
$ws.Range("E41").Value = '  +34.79%  '

# Row 42
$ws.Range("E42").Value = '  -0.82%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.90'
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.088.68'
$ws.Range("E44").Value = '  +0.60%  '

# Row 45
$ws.Range("E45").Value = '  -2.57%  '

# Row 46
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.60'
$ws.Range("E47").Value = '  +0.81%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.00'
$ws.Range("E48").Value = '  +2.38%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.801.80'
$ws.Range("E49").Value = '  +0.11%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.99'
$ws.Range("E50").Value = '  +7.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.66'
$ws.Range("E51").Value = '  -0.85%  '
